$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The F1/F2 crosses are being split into distinct a/b sub-crosses (F1 -> F1a/F1b,
# F2 -> F2a/F2b) in column A.
$ws.Range("A3").Value = "F1a"
$ws.Range("A4").Value = "F1b"
$ws.Range("A5").Value = "F2a"
$ws.Range("A6").Value = "F2b"
$ws.Range("A7").Value = "P2"

# Update the active selection to match the author's saved cursor position.
$ws.Range("F10").Select()
